$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47; this shifts existing rows 47-112 down to 48-113
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with data (copy of former row 47 template,
# with updated Fecha/Volumen/Precio columns reflecting the new record)
$ws.Cells.Item(47, 1).Value = 1
$ws.Cells.Item(47, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(47, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(47, 4).Value = "2023-09-06"
$ws.Cells.Item(47, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 5).Value = 15
$ws.Cells.Item(47, 6).Value = 100112012
$ws.Cells.Item(47, 7).Value = "Espinaca"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 285
$ws.Cells.Item(47, 11).Value = 2800
$ws.Cells.Item(47, 12).Value = 3000
$ws.Cells.Item(47, 13).Value = 2895
$ws.Cells.Item(47, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 965
$ws.Cells.Item(47, 17).Value = 3
$ws.Cells.Item(47, 18).Value = "Hortaliza"
